$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1 ("Data 1"): add two new supporting-material columns (W, X) ---
$ws1.Range("W1").Value = "https://www.ibm.com/it-it/topics/patch-management"
$ws1.Range("X1").Value = "https://cpl.thalesgroup.com/it/software-monetization/what-is-patch-management"

$ws1.Hyperlinks.Add($ws1.Range("W1"), "https://www.ibm.com/it-it/topics/patch-management")
$ws1.Hyperlinks.Add($ws1.Range("X1"), "https://cpl.thalesgroup.com/it/software-monetization/what-is-patch-management")

# --- Sheet 2 ("data as rows"): add the same two items as new rows (30, 31) ---
$ws2.Range("A30").Value = "https://www.ibm.com/it-it/topics/patch-management"
$ws2.Range("A31").Value = "https://cpl.thalesgroup.com/it/software-monetization/what-is-patch-management"

$ws2.Hyperlinks.Add($ws2.Range("A30"), "https://www.ibm.com/it-it/topics/patch-management")
$ws2.Hyperlinks.Add($ws2.Range("A31"), "https://cpl.thalesgroup.com/it/software-monetization/what-is-patch-management")

# --- Selections on each sheet following the edits ---
$ws1.Range("X1").Select()
$ws2.Range("A31").Select()

# --- Activate sheet 2 ("data as rows") last so it becomes the active tab ---
$ws2.Activate()
